# Metro budget exercise - add "Question 7" (7.3/7.4/7.5) lookups in columns F:I,
# mirroring the existing Question 3/4/5 exercises (columns A:D), but looking up
# by department name placed in column F and resolving the target column via
# MATCH against the header row, instead of relying on the raw table order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# ---------------------------------------------------------------------------
# Question 7.3 (VLOOKUP) - mirrors "Question 3" block at rows 54-61
# ---------------------------------------------------------------------------

# Section title
$ws.Range("F54").Value = "Question 7.3"
$ws.Range("A54").Copy()
$ws.Range("F54").PasteSpecial(-4122)

# Header row (Department / FY17_diff / FY18_diff / FY19_diff)
$ws.Range("F55").Value = "Department"
$ws.Range("G55").Value = "FY17_diff"
$ws.Range("H55").Value = "FY18_diff"
$ws.Range("I55").Value = "FY19_diff"
$ws.Range("A55:D55").Copy()
$ws.Range("F55").PasteSpecial(-4122)

# Re-write the original B56:D61 VLOOKUP formulas using a literal column index
# instead of MATCH(...) against the header.
$ws.Range("B56").Formula = "=VLOOKUP(`$A56, `$A`$2:`$P`$52, 4, FALSE)"
$ws.Range("C56").Formula = "=VLOOKUP(`$A56, `$A`$2:`$P`$52, 9, FALSE)"
$ws.Range("D56").Formula = "=VLOOKUP(`$A56, `$A`$2:`$P`$52, 14, FALSE)"
$ws.Range("B57:B61").Formula = "=VLOOKUP(`$A57, `$A`$2:`$P`$52, 4, FALSE)"
$ws.Range("C57:C61").Formula = "=VLOOKUP(`$A57, `$A`$2:`$P`$52, 9, FALSE)"
$ws.Range("D57:D61").Formula = "=VLOOKUP(`$A57, `$A`$2:`$P`$52, 14, FALSE)"

# New F:I department/value columns for rows 56-61 - department names mirror
# column A, values are VLOOKUP keyed off the new header row via MATCH.
$ws.Range("F56").Value = "Community Education Commission"
$ws.Range("F57").Value = "Community Oversight Board"
$ws.Range("F58").Value = "Election Commission"
$ws.Range("F59").Value = "Historical Commission"
$ws.Range("F60").Value = "Human Relations Commission"
$ws.Range("F61").Value = "Planning Commission"

$ws.Range("G56").Formula = "=VLOOKUP(`$A56, `$A`$2:`$P`$52, MATCH(G`$55, `$1:`$1, 0), FALSE)"
$ws.Range("H56:I56").Formula = "=VLOOKUP(`$A56, `$A`$2:`$P`$52, MATCH(H`$55, `$1:`$1, 0), FALSE)"
$ws.Range("G57:I61").Formula = "=VLOOKUP(`$A57, `$A`$2:`$P`$52, MATCH(G`$55, `$1:`$1, 0), FALSE)"

# ---------------------------------------------------------------------------
# Question 7.4 (XLOOKUP) - mirrors "Question 4" block at rows 63-70
# ---------------------------------------------------------------------------

$ws.Range("F63").Value = "Question 7.4"
$ws.Range("A63").Copy()
$ws.Range("F63").PasteSpecial(-4122)

$ws.Range("F64").Value = "Department"
$ws.Range("G64").Value = "FY17_diff"
$ws.Range("H64").Value = "FY18_diff"
$ws.Range("I64").Value = "FY19_diff"
$ws.Range("A64:D64").Copy()
$ws.Range("F64").PasteSpecial(-4122)

$ws.Range("F65").Value = "Community Education Commission"
$ws.Range("F66").Value = "Community Oversight Board"
$ws.Range("F67").Value = "Election Commission"
$ws.Range("F68").Value = "Historical Commission"
$ws.Range("F69").Value = "Human Relations Commission"
$ws.Range("F70").Value = "Planning Commission"

# Each cell gets its own (non-shared) single-cell array formula, matching
# what Excel produces when an XLOOKUP/INDEX dynamic-array formula is filled
# down one row at a time.
$ws.Range("G65").FormulaArray = "=_xlfn.XLOOKUP(`$F65,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("H65").FormulaArray = "=_xlfn.XLOOKUP(`$F65,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("I65").FormulaArray = "=_xlfn.XLOOKUP(`$F65,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$55,`$A`$1:`$P`$1,0)),FALSE)"

$ws.Range("G66").FormulaArray = "=_xlfn.XLOOKUP(`$F66,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("H66").FormulaArray = "=_xlfn.XLOOKUP(`$F66,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("I66").FormulaArray = "=_xlfn.XLOOKUP(`$F66,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$55,`$A`$1:`$P`$1,0)),FALSE)"

$ws.Range("G67").FormulaArray = "=_xlfn.XLOOKUP(`$F67,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("H67").FormulaArray = "=_xlfn.XLOOKUP(`$F67,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("I67").FormulaArray = "=_xlfn.XLOOKUP(`$F67,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$55,`$A`$1:`$P`$1,0)),FALSE)"

$ws.Range("G68").FormulaArray = "=_xlfn.XLOOKUP(`$F68,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("H68").FormulaArray = "=_xlfn.XLOOKUP(`$F68,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("I68").FormulaArray = "=_xlfn.XLOOKUP(`$F68,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$55,`$A`$1:`$P`$1,0)),FALSE)"

$ws.Range("G69").FormulaArray = "=_xlfn.XLOOKUP(`$F69,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("H69").FormulaArray = "=_xlfn.XLOOKUP(`$F69,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("I69").FormulaArray = "=_xlfn.XLOOKUP(`$F69,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$55,`$A`$1:`$P`$1,0)),FALSE)"

$ws.Range("G70").FormulaArray = "=_xlfn.XLOOKUP(`$F70,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("H70").FormulaArray = "=_xlfn.XLOOKUP(`$F70,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$55,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("I70").FormulaArray = "=_xlfn.XLOOKUP(`$F70,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$55,`$A`$1:`$P`$1,0)),FALSE)"

# Blank styled spacer row between the XLOOKUP block and Question 7.5 below it.
$ws.Range("G71:I71").Font.Bold = $true
$ws.Range("G71:I71").Font.Bold = $false

# ---------------------------------------------------------------------------
# Question 7.5 (INDEX/MATCH) - mirrors "Question 5" block at rows 72-79
# ---------------------------------------------------------------------------

$ws.Range("F72").Value = "Question 7.5"
$ws.Range("A72").Copy()
$ws.Range("F72").PasteSpecial(-4122)
$ws.Range("G72:I72").Font.Bold = $true
$ws.Range("G72:I72").Font.Bold = $false

$ws.Range("F73").Value = "Department"
$ws.Range("G73").Value = "FY17_diff"
$ws.Range("H73").Value = "FY18_diff"
$ws.Range("I73").Value = "FY19_diff"
$ws.Range("A73").Copy()
$ws.Range("F73").PasteSpecial(-4122)
$ws.Range("B55").Copy()
$ws.Range("G73:I73").PasteSpecial(-4122)

$ws.Range("F74").Value = "Community Education Commission"
$ws.Range("F75").Value = "Community Oversight Board"
$ws.Range("F76").Value = "Election Commission"
$ws.Range("F77").Value = "Historical Commission"
$ws.Range("F78").Value = "Human Relations Commission"
$ws.Range("F79").Value = "Planning Commission"

$ws.Range("G74").Formula = "=INDEX(`$A`$1:`$P`$52,MATCH(`$F74,`$A`$1:`$A`$52,0),MATCH(G`$73,`$A`$1:`$P`$1,0))"
$ws.Range("H74:I74").Formula = "=INDEX(`$A`$1:`$P`$52,MATCH(`$F74,`$A`$1:`$A`$52,0),MATCH(H`$73,`$A`$1:`$P`$1,0))"
$ws.Range("G75:I79").Formula = "=INDEX(`$A`$1:`$P`$52,MATCH(`$F75,`$A`$1:`$A`$52,0),MATCH(G`$73,`$A`$1:`$P`$1,0))"

# Re-apply the value-only style ("10"-equivalent) look to the new lookup
# columns so the block reads as a clean table (no-op visual touch forces the
# cells to materialise with the workbook's default numeric style).
$ws.Range("G74:I79").Font.Bold = $true
$ws.Range("G74:I79").Font.Bold = $false

# ---------------------------------------------------------------------------
# View state: scroll position / active selection left where the author's
# last edit was (inside the new Question 7.4 XLOOKUP block).
# ---------------------------------------------------------------------------
$ws.Range("K69").Select()
